$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.672.71'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '3.486.14'
$ws.Range('E3').Value = '  +5.75%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '261.47'
$ws.Range('E5').Value = '  +2.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '674.60'
$ws.Range('E6').Value = '  +8.76%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.55'
$ws.Range('E7').Value = '  +9.97%  '
$ws.Range('E8').Value = '  +16.99%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.10'
$ws.Range('E9').Value = '  +23.08%  '
$ws.Range('E10').Value = '  -0.15%  '
$ws.Range('D11').Value = '3.483.78'
$ws.Range('E11').Value = '  +5.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.223'
$ws.Range('E12').Value = '  +12.52%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '42.67'
$ws.Range('E13').Value = '  +11.28%  '
$ws.Range('E14').Value = '  +11.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.21'
$ws.Range('E15').Value = '  +14.29%  '
$ws.Range('D16').Value = '98.483.59'
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('D17').Value = '4.128.55'
$ws.Range('E17').Value = '  +5.44%  '
$ws.Range('E18').Value = '  +34.19%  '
$ws.Range('D19').Value = '3.476.22'
$ws.Range('E19').Value = '  +5.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.60'
$ws.Range('E20').Value = '  +16.97%  '
$ws.Range('B21').Value = 'SuiNetwork'
$ws.Range('C21').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.64'
$ws.Range('E21').Value = '  +3.78%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '538.80'
$ws.Range('E22').Value = '  +13.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.84'
$ws.Range('E23').Value = '  +15.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000222'
$ws.Range('E24').Value = '  +8.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.457'
$ws.Range('E25').Value = '  +57.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.42'
$ws.Range('E26').Value = '  +15.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '102.93'
$ws.Range('E27').Value = '  +17.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.99'
$ws.Range('E28').Value = '  +9.97%  '
$ws.Range('D29').Value = '3.654.17'
$ws.Range('E29').Value = '  +5.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.151'
$ws.Range('E30').Value = '  +16.40%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.51'
$ws.Range('E31').Value = '  +18.09%  '
$ws.Range('B32').Value = 'Cronos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.200'
$ws.Range('E32').Value = '  +7.66%  '
$ws.Range('E33').Value = '  +0.15%  '
$ws.Range('B34').Value = 'PolygonEcosystemToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.583'
$ws.Range('E34').Value = '  +29.18%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '30.84'
$ws.Range('E35').Value = '  +12.39%  '
$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.23'
$ws.Range('E37').Value = '  +16.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.01'
$ws.Range('E38').Value = '  +12.36%  '
$ws.Range('E39').Value = '  +10.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '542.86'
$ws.Range('E40').Value = '  +11.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.44'
$ws.Range('E41').Value = '  +16.66%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.875'
$ws.Range('E43').Value = '  +9.76%  '
$ws.Range('E44').Value = '  +35.54%  '
$ws.Range('E45').Value = '  +12.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.77'
$ws.Range('E46').Value = '  +3.42%  '
$ws.Range('B47').Value = 'Cosmos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.27'
$ws.Range('E47').Value = '  +18.36%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.14'
$ws.Range('E48').Value = '  +12.93%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.59'
$ws.Range('E50').Value = '  +18.70%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.28'
$ws.Range('E51').Value = '  +15.14%  '
